$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 29.125
$ws.Range("I5").Value = 30.428572
$ws.Range("J5").Value = 20
$ws.Range("K5").Value = 30.428572
$ws.Range("L5").Value = 20
$ws.Range("M5").Value = 84.571428
$ws.Range("N5").Value = -250
# Row 10
$ws.Range("H10").Value = 3500
$ws.Range("I10").Value = 2000
$ws.Range("J10").Value = 5000
$ws.Range("K10").Value = 2000
$ws.Range("L10").Value = 5000
$ws.Range("M10").Value = -1707
$ws.Range("N10").Value = -5586
# Row 18
$ws.Range("H18").Value = 3232
$ws.Range("I18").Value = 2680
$ws.Range("J18").Value = 5992
$ws.Range("K18").Value = 2680
$ws.Range("L18").Value = 5992
$ws.Range("M18").Value = -2396
$ws.Range("N18").Value = -6560
# Row 100
$ws.Range("H100").Value = 7748
$ws.Range("I100").Value = 8489
$ws.Range("J100").Value = 7599.8
$ws.Range("K100").Value = 8489
$ws.Range("L100").Value = 7599.8
$ws.Range("M100").Value = -7948
$ws.Range("N100").Value = -8681.799999999999
# Row 101
$ws.Range("H101").Value = 1269.5555
$ws.Range("I101").Value = 248.66667
$ws.Range("K101").Value = 746.00001
$ws.Range("M101").Value = 875.99999
# Row 107
$ws.Range("H107").Value = 1862.4736
$ws.Range("J107").Value = 1435
$ws.Range("L107").Value = 1435
$ws.Range("N107").Value = -5275
# Row 111
$ws.Range("H111").Value = 858.5
$ws.Range("I111").Value = 962.6667
$ws.Range("J111").Value = 754.3333
$ws.Range("K111").Value = 2888.0001
$ws.Range("L111").Value = 2262.9999
$ws.Range("M111").Value = 178.9998999999998
$ws.Range("N111").Value = -8396.999899999999
# Row 113
$ws.Range("H113").Value = 3894.25
$ws.Range("I113").Value = 3897.5217
$ws.Range("J113").Value = 3879.2
$ws.Range("K113").Value = 3897.5217
$ws.Range("L113").Value = 3879.2
$ws.Range("M113").Value = -643.5216999999998
$ws.Range("N113").Value = -10387.2
# Row 116
$ws.Range("H116").Value = 10297.5
$ws.Range("J116").Value = 10276.4
$ws.Range("L116").Value = 10276.4
$ws.Range("N116").Value = -17160.4
# Row 137
$ws.Range("H137").Value = 750910.2
$ws.Range("I137").Value = 1083.625
$ws.Range("J137").Value = 1673773.8
$ws.Range("K137").Value = 3250.875
$ws.Range("L137").Value = 5021321.4
$ws.Range("M137").Value = -700.875
$ws.Range("N137").Value = -5026421.4

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 3
$ws.Range("H3").Value = 5583.6665
$ws.Range("I3").Value = 4300.4
$ws.Range("K3").Value = 4300.4
$ws.Range("M3").Value = -4185.4
# Row 32
$ws.Range("H32").Value = 143832.25
$ws.Range("I32").Value = 143832.25
$ws.Range("K32").Value = 143832.25
$ws.Range("M32").Value = -143545.25
# Row 74
$ws.Range("H74").Value = 1974.6604
$ws.Range("I74").Value = 1651.3182
$ws.Range("J74").Value = 3555.4443
$ws.Range("K74").Value = 1651.3182
$ws.Range("L74").Value = 3555.4443
$ws.Range("M74").Value = -777.3181999999999
$ws.Range("N74").Value = -5303.4443
# Row 77
$ws.Range("H77").Value = 1974.6604
$ws.Range("I77").Value = 1651.3182
$ws.Range("J77").Value = 3555.4443
$ws.Range("K77").Value = 8256.591
$ws.Range("L77").Value = 17777.2215
$ws.Range("M77").Value = -3888.591
$ws.Range("N77").Value = -26513.2215
# Row 122
$ws.Range("H122").Value = 55605556
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
# Row 132
$ws.Range("H132").Value = 1545.52
$ws.Range("I132").Value = 1476.4286
$ws.Range("K132").Value = 4429.2858
$ws.Range("M132").Value = -1899.2858

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Range("H5").Value = 497.9375
$ws.Range("J5").Value = 1016.8333
$ws.Range("L5").Value = 1016.8333
$ws.Range("N5").Value = -1242.8333
# Row 107
$ws.Range("H107").Value = 4843.1177
$ws.Range("I107").Value = 3029.5
$ws.Range("K107").Value = 3029.5
$ws.Range("M107").Value = -1109.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 3236
$ws.Range("J16").Value = 3930.3333
$ws.Range("L16").Value = 3930.3333
$ws.Range("N16").Value = -4504.3333
# Row 19
$ws.Range("H19").Value = 549.3461
$ws.Range("I19").Value = 556.65
$ws.Range("J19").Value = 525
$ws.Range("K19").Value = 556.65
$ws.Range("L19").Value = 525
$ws.Range("M19").Value = -386.65
$ws.Range("N19").Value = -865
# Row 22
$ws.Range("H22").Value = 1393.25
$ws.Range("I22").Value = 575
$ws.Range("J22").Value = 1666
$ws.Range("K22").Value = 575
$ws.Range("L22").Value = 1666
$ws.Range("M22").Value = -225
$ws.Range("N22").Value = -2366
# Row 24
$ws.Range("H24").Value = 549.3461
$ws.Range("I24").Value = 556.65
$ws.Range("J24").Value = 525
$ws.Range("K24").Value = 556.65
$ws.Range("L24").Value = 525
$ws.Range("M24").Value = -386.65
$ws.Range("N24").Value = -865
# Row 31
$ws.Range("H31").Value = 3084.7693
$ws.Range("I31").Value = 1194.7693
$ws.Range("J31").Value = 4974.769
$ws.Range("K31").Value = 1194.7693
$ws.Range("L31").Value = 4974.769
$ws.Range("M31").Value = -899.7692999999999
$ws.Range("N31").Value = -5564.769
# Row 34
$ws.Range("H34").Value = 3084.7693
$ws.Range("I34").Value = 1194.7693
$ws.Range("J34").Value = 4974.769
$ws.Range("K34").Value = 1194.7693
$ws.Range("L34").Value = 4974.769
$ws.Range("M34").Value = -992.7692999999999
$ws.Range("N34").Value = -5378.769
# Row 41
$ws.Range("H41").Value = 29391.857
$ws.Range("J41").Value = 51246.25
$ws.Range("L41").Value = 51246.25
$ws.Range("N41").Value = -52102.25
# Row 64
$ws.Range("H64").Value = 74999.5
$ws.Range("J64").Value = 74999.5
$ws.Range("L64").Value = 74999.5
$ws.Range("N64").Value = -75495.5
# Row 67
$ws.Range("H67").Value = 74999.5
$ws.Range("J67").Value = 74999.5
$ws.Range("L67").Value = 74999.5
$ws.Range("N67").Value = -76715.5
# Row 86
$ws.Range("H86").Value = 6982.5386
$ws.Range("I86").Value = 5876.1113
$ws.Range("K86").Value = 5876.1113
$ws.Range("M86").Value = -4753.1113
# Row 89
$ws.Range("H89").Value = 6982.5386
$ws.Range("I89").Value = 5876.1113
$ws.Range("K89").Value = 29380.5565
$ws.Range("M89").Value = -23764.5565
# Row 105
$ws.Range("H105").Value = 4262.3
$ws.Range("I105").Value = 3139.3845
$ws.Range("J105").Value = 6347.7144
$ws.Range("K105").Value = 3139.3845
$ws.Range("L105").Value = 6347.7144
$ws.Range("M105").Value = -1392.3845
$ws.Range("N105").Value = -9841.714400000001
# Row 107
$ws.Range("H107").Value = 437.82758
$ws.Range("I107").Value = 302.33334
$ws.Range("K107").Value = 302.33334
$ws.Range("M107").Value = 1617.66666
# Row 113
$ws.Range("H113").Value = 3236
$ws.Range("J113").Value = 3930.3333
$ws.Range("L113").Value = 3930.3333
$ws.Range("N113").Value = -8270.3333
# Row 122
$ws.Range("H122").Value = 15874398
$ws.Range("I122").Value = 1841.091
$ws.Range("J122").Value = 33334210
$ws.Range("K122").Value = 5523.272999999999
$ws.Range("L122").Value = 100002630
$ws.Range("M122").Value = -3073.272999999999
$ws.Range("N122").Value = -100007530
# Row 141
$ws.Range("H141").Value = 204499.55
$ws.Range("J141").Value = 204499.55
$ws.Range("L141").Value = 204499.55
$ws.Range("N141").Value = -214859.55

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 537.8461
$ws.Range("I14").Value = 537.8461
$ws.Range("K14").Value = 1613.5383
$ws.Range("M14").Value = -1440.5383
# Row 113
$ws.Range("H113").Value = 1179.2307
$ws.Range("I113").Value = 1774.6
$ws.Range("J113").Value = 807.125
$ws.Range("K113").Value = 5323.799999999999
$ws.Range("L113").Value = 2421.375
$ws.Range("M113").Value = -3153.799999999999
$ws.Range("N113").Value = -6761.375

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 22
$ws.Range("H22").Value = 916.6667
$ws.Range("I22").Value = 916.6667
$ws.Range("K22").Value = 916.6667
$ws.Range("M22").Value = -387.6667
# Row 80
$ws.Range("H80").Value = 76928810
$ws.Range("I80").Value = 333336000
$ws.Range("K80").Value = 333336000
$ws.Range("M80").Value = -333335002
# Row 83
$ws.Range("H83").Value = 76928810
$ws.Range("I83").Value = 333336000
$ws.Range("K83").Value = 1666680000
$ws.Range("M83").Value = -1666675008
# Row 102
$ws.Range("H102").Value = 1082.2
$ws.Range("I102").Value = 775.6957
$ws.Range("K102").Value = 775.6957
$ws.Range("M102").Value = 846.3043
# Row 107
$ws.Range("H107").Value = 1905
$ws.Range("I107").Value = 494.4
$ws.Range("J107").Value = 2912.5715
$ws.Range("K107").Value = 494.4
$ws.Range("L107").Value = 2912.5715
$ws.Range("M107").Value = 1425.6
$ws.Range("N107").Value = -6752.5715
# Row 122
$ws.Range("H122").Value = 8107.7896
$ws.Range("J122").Value = 5499.75
$ws.Range("L122").Value = 16499.25
$ws.Range("N122").Value = -21399.25
# Row 126
$ws.Range("H126").Value = 2368.625
$ws.Range("I126").Value = 1813.5
$ws.Range("K126").Value = 5440.5
$ws.Range("M126").Value = -2970.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 1683.5834
$ws.Range("I55").Value = 1245.1111
$ws.Range("K55").Value = 1245.1111
$ws.Range("M55").Value = -1072.1111
# Row 70
$ws.Range("H70").Value = 333368320
$ws.Range("J70").Value = 333368320
$ws.Range("L70").Value = 333368320
$ws.Range("N70").Value = -333368860
# Row 73
$ws.Range("H73").Value = 333368320
$ws.Range("J73").Value = 333368320
$ws.Range("L73").Value = 333368320
$ws.Range("N73").Value = -333370192

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 24
$ws.Range("H24").Value = 15500
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 15500
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 15500
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -15960
# Row 122
$ws.Range("H122").Value = 531490.4399999999
$ws.Range("I122").Value = 1151593
$ws.Range("K122").Value = 3454779
$ws.Range("M122").Value = -3452329
# Row 132
$ws.Range("H132").Value = 33336248
$ws.Range("I132").Value = 33336248
$ws.Range("K132").Value = 100008744
$ws.Range("M132").Value = -100006214
